$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1806.5333
$ws.Range("I15").Value = 1806.5333
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5419.5999
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5250.5999

$ws.Range("H18").Value = 3992.5
$ws.Range("I18").Value = 3992.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3992.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3708.5
$ws.Range("N18").ClearContents()

$ws.Range("H116").Value = 4849.5
$ws.Range("I116").Value = 4849
$ws.Range("J116").Value = 4850
$ws.Range("K116").Value = 4849
$ws.Range("L116").Value = 4850
$ws.Range("M116").Value = -1407
$ws.Range("N116").Value = -11734

$ws.Range("H135").Value = 3232.75
$ws.Range("I135").Value = 3904.3333
$ws.Range("J135").Value = 1218
$ws.Range("K135").Value = 35138.9997
$ws.Range("L135").Value = 10962
$ws.Range("M135").Value = -32603.9997

$ws.Range("H137").Value = 6107.2
$ws.Range("I137").Value = 12638.8
$ws.Range("J137").Value = 4474.3
$ws.Range("K137").Value = 37916.39999999999
$ws.Range("L137").Value = 13422.9
$ws.Range("M137").Value = -35366.39999999999

$ws.Range("H138").Value = 3843.6904
$ws.Range("I138").Value = 2644.8462
$ws.Range("J138").Value = 4381.1035
$ws.Range("K138").Value = 7934.5386
$ws.Range("L138").Value = 13143.3105
$ws.Range("M138").Value = -2794.5386
$ws.Range("N138").Value = -23423.3105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5719.8335
$ws.Range("I32").Value = 4648.9546
$ws.Range("J32").Value = 17499.5
$ws.Range("K32").Value = 4648.9546
$ws.Range("L32").Value = 17499.5
$ws.Range("M32").Value = -4361.9546

$ws.Range("H45").Value = 2158.9412
$ws.Range("I45").Value = 2039.8889
$ws.Range("J45").Value = 2292.875
$ws.Range("K45").Value = 2039.8889
$ws.Range("L45").Value = 2292.875
$ws.Range("M45").Value = -1662.8889
$ws.Range("N45").Value = -3046.875

$ws.Range("H61").Value = 1875.2858
$ws.Range("I61").Value = 1712.0385
$ws.Range("J61").Value = 3997.5
$ws.Range("K61").Value = 1712.0385
$ws.Range("L61").Value = 3997.5
$ws.Range("M61").Value = -1500.0385
$ws.Range("N61").Value = -4421.5

$ws.Range("H122").Value = 9666.666999999999
$ws.Range("I122").Value = 9666.666999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29000.001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -26550.001

$ws.Range("H132").Value = 2815.2173
$ws.Range("I132").Value = 1848.2941
$ws.Range("J132").Value = 5554.8335
$ws.Range("K132").Value = 5544.8823
$ws.Range("L132").Value = 16664.5005
$ws.Range("M132").Value = -3014.8823

$ws.Range("H136").Value = 1875.2858
$ws.Range("I136").Value = 1712.0385
$ws.Range("J136").Value = 3997.5
$ws.Range("K136").Value = 5136.1155
$ws.Range("L136").Value = 11992.5
$ws.Range("M136").Value = -2586.1155
$ws.Range("N136").Value = -17092.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1171.7142
$ws.Range("I86").Value = 1171.7142
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1171.7142
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -48.71419999999989

$ws.Range("H89").Value = 1171.7142
$ws.Range("I89").Value = 1171.7142
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5858.571
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -242.5709999999999

$ws.Range("H99").Value = 803.3
$ws.Range("I99").Value = 814
$ws.Range("J99").Value = 707
$ws.Range("K99").Value = 814
$ws.Range("L99").Value = 707
$ws.Range("M99").Value = 684
$ws.Range("N99").Value = -3703

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6253.857
$ws.Range("I86").Value = 6233.5557
$ws.Range("J86").Value = 6290.4
$ws.Range("K86").Value = 6233.5557
$ws.Range("L86").Value = 6290.4
$ws.Range("M86").Value = -5110.5557
$ws.Range("N86").Value = -8536.4

$ws.Range("H89").Value = 6253.857
$ws.Range("I89").Value = 6233.5557
$ws.Range("J89").Value = 6290.4
$ws.Range("K89").Value = 31167.7785
$ws.Range("L89").Value = 31452
$ws.Range("M89").Value = -25551.7785
$ws.Range("N89").Value = -42684

$ws.Range("H99").Value = 1999.8334
$ws.Range("I99").Value = 2159.8
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 2159.8
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = -661.8000000000002
$ws.Range("N99").Value = -4196

$ws.Range("H126").Value = 1999.8334
$ws.Range("I126").Value = 2159.8
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 6479.400000000001
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -4009.400000000001
$ws.Range("N126").Value = -8540

$ws.Range("H132").Value = 4257.8
$ws.Range("I132").Value = 3541.889
$ws.Range("J132").Value = 5331.6665
$ws.Range("K132").Value = 10625.667
$ws.Range("L132").Value = 15994.9995
$ws.Range("M132").Value = -8095.667000000001
$ws.Range("N132").Value = -21054.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1941.6666
$ws.Range("I107").Value = 2704.2
$ws.Range("J107").Value = 988.5
$ws.Range("K107").Value = 8112.599999999999
$ws.Range("L107").Value = 2965.5
$ws.Range("M107").Value = -6192.599999999999
$ws.Range("N107").Value = -6805.5

$ws.Range("H113").Value = 688.4
$ws.Range("I113").Value = 399.4
$ws.Range("J113").Value = 977.4
$ws.Range("K113").Value = 1198.2
$ws.Range("L113").Value = 2932.2
$ws.Range("M113").Value = 971.8000000000002
$ws.Range("N113").Value = -7272.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 974.2857
$ws.Range("I97").Value = 803.3333
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 803.3333
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -307.3333

$ws.Range("H102").Value = 4098.3335
$ws.Range("I102").Value = 4098.3335
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4098.3335
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2476.3335

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 990.4
$ws.Range("I122").Value = 989.55554
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 2968.66662
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -518.66662

$ws.Range("H132").Value = 3111.7727
$ws.Range("I132").Value = 2655.9473
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 7967.841899999999
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -5437.841899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3969.2307
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3969.2307
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3969.2307
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -4559.2307

$ws.Range("H27").Value = 3969.2307
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3969.2307
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3969.2307
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4183.2307

$ws.Range("H55").Value = 500.83334
$ws.Range("I55").Value = 561.6667
$ws.Range("J55").Value = 440
$ws.Range("K55").Value = 561.6667
$ws.Range("L55").Value = 440
$ws.Range("M55").Value = -388.6667
$ws.Range("N55").Value = -786

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 44643.6
$ws.Range("I70").Value = 39999.5
$ws.Range("J70").Value = 45358.08
$ws.Range("K70").Value = 39999.5
$ws.Range("L70").Value = 45358.08
$ws.Range("M70").Value = -39684.5
$ws.Range("N70").Value = -45988.08

$ws.Range("H73").Value = 44643.6
$ws.Range("I73").Value = 39999.5
$ws.Range("J73").Value = 45358.08
$ws.Range("K73").Value = 39999.5
$ws.Range("L73").Value = 45358.08
$ws.Range("M73").Value = -38907.5
$ws.Range("N73").Value = -47542.08

$ws.Range("H107").Value = 313
$ws.Range("I107").Value = 282.16666
$ws.Range("J107").Value = 498
$ws.Range("K107").Value = 846.4999799999999
$ws.Range("L107").Value = 1494
$ws.Range("M107").Value = 1073.50002

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
